$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Environmental_variables")
$ws2 = $wb.Worksheets.Item("Monthly_temp_and_fog")

# Sheet1 header renames: *_aug -> *_jul
$ws1.Range("F1").Value = "freq_jul"
$ws1.Range("I1").Value = "temp_jul"
$ws1.Range("K1").Value = "RH_jul"
$ws1.Range("M1").Value = "percip_jul"

# Sheet1 data updates: replace August values with July values
$ws1.Range("F2").Value = 0.049212597310543102
$ws1.Range("I2").Value = 17.535275537634401
$ws1.Range("K2").Value = 88.589669578853005
$ws1.Range("M2").Value = 244.10330149999999
$ws1.Range("F3").Value = 0.039215687662362997
$ws1.Range("I3").Value = 17.358539426523301
$ws1.Range("K3").Value = 83.344968637993006
$ws1.Range("M3").Value = 247.93342340000001
$ws1.Range("F4").Value = 0.039215687662362997
$ws1.Range("I4").Value = 17.3162939068099
$ws1.Range("K4").Value = 91.003519937275996
$ws1.Range("M4").Value = 247.93342340000001
$ws1.Range("F5").Value = 0.04296875
$ws1.Range("I5").Value = 17.478683691756199
$ws1.Range("K5").Value = 89.187310035842401
$ws1.Range("M5").Value = 246.9084618
$ws1.Range("F6").Value = 0.04296875
$ws1.Range("I6").Value = 17.609603158602098
$ws1.Range("K6").Value = 89.028451388888897
$ws1.Range("M6").Value = 246.9084618
$ws1.Range("F7").Value = 0.047058824449777603
$ws1.Range("I7").Value = 17.740522625448001
$ws1.Range("K7").Value = 88.869592741935506
$ws1.Range("M7").Value = 246.9084618
$ws1.Range("F8").Value = 0.0459770113229752
$ws1.Range("I8").Value = 17.280460125447998
$ws1.Range("K8").Value = 91.875328629032396
$ws1.Range("M8").Value = 248.81373139999999
$ws1.Range("F9").Value = 0.0555555559694767
$ws1.Range("I9").Value = 17.171097222222201
$ws1.Range("K9").Value = 86.818875448028706
$ws1.Range("M9").Value = 248.81373139999999
$ws1.Range("F10").Value = 0.057361375540494898
$ws1.Range("I10").Value = 17.162118951612999
$ws1.Range("K10").Value = 93.640293682795701
$ws1.Range("M10").Value = 248.81373139999999
$ws1.Range("F11").Value = 0.083419999999999994
$ws1.Range("I11").Value = 16.8830613799283
$ws1.Range("K11").Value = 91.076397177419395
$ws1.Range("M11").Value = 228.8770079
$ws1.Range("F12").Value = 0.079150579869747203
$ws1.Range("I12").Value = 17.306602150537699
$ws1.Range("K12").Value = 90.137686379928397
$ws1.Range("M12").Value = 230.33842139999999
$ws1.Range("F13").Value = 0.083419999999999994
$ws1.Range("I13").Value = 17.094831765233
$ws1.Range("K13").Value = 90.607041778673903
$ws1.Range("M13").Value = 228.8770079
$ws1.Range("F14").Value = 0.10150375962257401
$ws1.Range("I14").Value = 17.163833557347601
$ws1.Range("K14").Value = 91.776940860215007
$ws1.Range("M14").Value = 248.4360777
$ws1.Range("F15").Value = 0.10150375962257401
$ws1.Range("I15").Value = 17.224108646953301
$ws1.Range("K15").Value = 92.393137544802798
$ws1.Range("M15").Value = 247.6530128
$ws1.Range("F16").Value = 0.073170728981494904
$ws1.Range("I16").Value = 16.8538508064517
$ws1.Range("K16").Value = 86.588230286738394
$ws1.Range("M16").Value = 250.6924147
$ws1.Range("F17").Value = 0.12098298966884601
$ws1.Range("I17").Value = 17.767870743727599
$ws1.Range("K17").Value = 92.678935259856601
$ws1.Range("M17").Value = 250.53075440000001
$ws1.Range("F18").Value = 0.11641221493482599
$ws1.Range("I18").Value = 17.754681451612999
$ws1.Range("K18").Value = 86.119648297490798
$ws1.Range("M18").Value = 250.53075440000001
$ws1.Range("F19").Value = 0.13307984173297899
$ws1.Range("I19").Value = 17.378605734767
$ws1.Range("K19").Value = 93.943468637992893
$ws1.Range("M19").Value = 250.61170659999999
$ws1.Range("F20").Value = 0.156673118472099
$ws1.Range("I20").Value = 18.060225358423001
$ws1.Range("K20").Value = 91.119783154121905
$ws1.Range("M20").Value = 226.62886
$ws1.Range("F21").Value = 0.130350187420845
$ws1.Range("I21").Value = 18.070809587813699
$ws1.Range("K21").Value = 89.785393593189994
$ws1.Range("M21").Value = 226.62886
$ws1.Range("F22").Value = 0.130350187420845
$ws1.Range("I22").Value = 17.862373879928299
$ws1.Range("K22").Value = 85.533505824372398
$ws1.Range("M22").Value = 226.62886
$ws1.Range("F23").Value = 0.173076927661896
$ws1.Range("I23").Value = 17.9100477150538
$ws1.Range("K23").Value = 90.713176747311906
$ws1.Range("M23").Value = 213.84870939999999
$ws1.Range("F24").Value = 0.14779271185398099
$ws1.Range("I24").Value = 18.222874103942701
$ws1.Range("K24").Value = 89.807600806451603
$ws1.Range("M24").Value = 210.75810250000001
$ws1.Range("F25").Value = 0.15296366810798601
$ws1.Range("I25").Value = 18.027835125448
$ws1.Range("K25").Value = 83.273996415770796
$ws1.Range("M25").Value = 210.75810250000001
$ws1.Range("F26").Value = 0.194931775331497
$ws1.Range("I26").Value = 18.040794130824398
$ws1.Range("K26").Value = 93.820918682795707
$ws1.Range("M26").Value = 256.93622069999998
$ws1.Range("F27").Value = 0.184675827622414
$ws1.Range("I27").Value = 17.8575112007168
$ws1.Range("K27").Value = 88.902224462365098
$ws1.Range("M27").Value = 256.93622069999998
$ws1.Range("F28").Value = 0.18834951519966101
$ws1.Range("I28").Value = 18.005357302867399
$ws1.Range("K28").Value = 91.957672043010803
$ws1.Range("M28").Value = 260.69217930000002

# Sheet2 header updates
$ws2.Range("K1").Value = "temp_may"
$ws2.Range("L1").Value = "temp_jun"
$ws2.Range("M1").Value = "temp_jul"
$ws2.Range("N1").Value = "fog_aug"
$ws2.Range("O1").Value = "fog_sep"
$ws2.Range("P1").Value = "fog_oct"
$ws2.Range("Q1").Value = "fog_nov"
$ws2.Range("R1").Value = "fog_dec"
$ws2.Range("S1").Value = "fog_jan"
$ws2.Range("T1").Value = "fog_feb"
$ws2.Range("U1").Value = "fog_mar"
$ws2.Range("V1").Value = "fog_apr"
$ws2.Range("W1").Value = "fog_may"
$ws2.Range("X1").Value = "fog_jun"
$ws2.Range("Y1").Value = "fog_jul"

# Sheet2 data: move old freq_* (K..S) values into new fog_* (N..V) positions
$__tmp = @{}
$__tmp["K"] = $ws2.Range("K2").Value2
$__tmp["L"] = $ws2.Range("L2").Value2
$__tmp["M"] = $ws2.Range("M2").Value2
$__tmp["N"] = $ws2.Range("N2").Value2
$__tmp["O"] = $ws2.Range("O2").Value2
$__tmp["P"] = $ws2.Range("P2").Value2
$__tmp["Q"] = $ws2.Range("Q2").Value2
$__tmp["R"] = $ws2.Range("R2").Value2
$__tmp["S"] = $ws2.Range("S2").Value2
$ws2.Range("N2").Value = $__tmp["K"]
$ws2.Range("O2").Value = $__tmp["L"]
$ws2.Range("P2").Value = $__tmp["M"]
$ws2.Range("Q2").Value = $__tmp["N"]
$ws2.Range("R2").Value = $__tmp["O"]
$ws2.Range("S2").Value = $__tmp["P"]
$ws2.Range("T2").Value = $__tmp["Q"]
$ws2.Range("U2").Value = $__tmp["R"]
$ws2.Range("V2").Value = $__tmp["S"]
$__tmp = @{}
$__tmp["K"] = $ws2.Range("K3").Value2
$__tmp["L"] = $ws2.Range("L3").Value2
$__tmp["M"] = $ws2.Range("M3").Value2
$__tmp["N"] = $ws2.Range("N3").Value2
$__tmp["O"] = $ws2.Range("O3").Value2
$__tmp["P"] = $ws2.Range("P3").Value2
$__tmp["Q"] = $ws2.Range("Q3").Value2
$__tmp["R"] = $ws2.Range("R3").Value2
$__tmp["S"] = $ws2.Range("S3").Value2
$ws2.Range("N3").Value = $__tmp["K"]
$ws2.Range("O3").Value = $__tmp["L"]
$ws2.Range("P3").Value = $__tmp["M"]
$ws2.Range("Q3").Value = $__tmp["N"]
$ws2.Range("R3").Value = $__tmp["O"]
$ws2.Range("S3").Value = $__tmp["P"]
$ws2.Range("T3").Value = $__tmp["Q"]
$ws2.Range("U3").Value = $__tmp["R"]
$ws2.Range("V3").Value = $__tmp["S"]
$__tmp = @{}
$__tmp["K"] = $ws2.Range("K4").Value2
$__tmp["L"] = $ws2.Range("L4").Value2
$__tmp["M"] = $ws2.Range("M4").Value2
$__tmp["N"] = $ws2.Range("N4").Value2
$__tmp["O"] = $ws2.Range("O4").Value2
$__tmp["P"] = $ws2.Range("P4").Value2
$__tmp["Q"] = $ws2.Range("Q4").Value2
$__tmp["R"] = $ws2.Range("R4").Value2
$__tmp["S"] = $ws2.Range("S4").Value2
$ws2.Range("N4").Value = $__tmp["K"]
$ws2.Range("O4").Value = $__tmp["L"]
$ws2.Range("P4").Value = $__tmp["M"]
$ws2.Range("Q4").Value = $__tmp["N"]
$ws2.Range("R4").Value = $__tmp["O"]
$ws2.Range("S4").Value = $__tmp["P"]
$ws2.Range("T4").Value = $__tmp["Q"]
$ws2.Range("U4").Value = $__tmp["R"]
$ws2.Range("V4").Value = $__tmp["S"]
$__tmp = @{}
$__tmp["K"] = $ws2.Range("K5").Value2
$__tmp["L"] = $ws2.Range("L5").Value2
$__tmp["M"] = $ws2.Range("M5").Value2
$__tmp["N"] = $ws2.Range("N5").Value2
$__tmp["O"] = $ws2.Range("O5").Value2
$__tmp["P"] = $ws2.Range("P5").Value2
$__tmp["Q"] = $ws2.Range("Q5").Value2
$__tmp["R"] = $ws2.Range("R5").Value2
$__tmp["S"] = $ws2.Range("S5").Value2
$ws2.Range("N5").Value = $__tmp["K"]
$ws2.Range("O5").Value = $__tmp["L"]
$ws2.Range("P5").Value = $__tmp["M"]
$ws2.Range("Q5").Value = $__tmp["N"]
$ws2.Range("R5").Value = $__tmp["O"]
$ws2.Range("S5").Value = $__tmp["P"]
$ws2.Range("T5").Value = $__tmp["Q"]
$ws2.Range("U5").Value = $__tmp["R"]
$ws2.Range("V5").Value = $__tmp["S"]
$__tmp = @{}
$__tmp["K"] = $ws2.Range("K6").Value2
$__tmp["L"] = $ws2.Range("L6").Value2
$__tmp["M"] = $ws2.Range("M6").Value2
$__tmp["N"] = $ws2.Range("N6").Value2
$__tmp["O"] = $ws2.Range("O6").Value2
$__tmp["P"] = $ws2.Range("P6").Value2
$__tmp["Q"] = $ws2.Range("Q6").Value2
$__tmp["R"] = $ws2.Range("R6").Value2
$__tmp["S"] = $ws2.Range("S6").Value2
$ws2.Range("N6").Value = $__tmp["K"]
$ws2.Range("O6").Value = $__tmp["L"]
$ws2.Range("P6").Value = $__tmp["M"]
$ws2.Range("Q6").Value = $__tmp["N"]
$ws2.Range("R6").Value = $__tmp["O"]
$ws2.Range("S6").Value = $__tmp["P"]
$ws2.Range("T6").Value = $__tmp["Q"]
$ws2.Range("U6").Value = $__tmp["R"]
$ws2.Range("V6").Value = $__tmp["S"]
$__tmp = @{}
$__tmp["K"] = $ws2.Range("K7").Value2
$__tmp["L"] = $ws2.Range("L7").Value2
$__tmp["M"] = $ws2.Range("M7").Value2
$__tmp["N"] = $ws2.Range("N7").Value2
$__tmp["O"] = $ws2.Range("O7").Value2
$__tmp["P"] = $ws2.Range("P7").Value2
$__tmp["Q"] = $ws2.Range("Q7").Value2
$__tmp["R"] = $ws2.Range("R7").Value2
$__tmp["S"] = $ws2.Range("S7").Value2
$ws2.Range("N7").Value = $__tmp["K"]
$ws2.Range("O7").Value = $__tmp["L"]
$ws2.Range("P7").Value = $__tmp["M"]
$ws2.Range("Q7").Value = $__tmp["N"]
$ws2.Range("R7").Value = $__tmp["O"]
$ws2.Range("S7").Value = $__tmp["P"]
$ws2.Range("T7").Value = $__tmp["Q"]
$ws2.Range("U7").Value = $__tmp["R"]
$ws2.Range("V7").Value = $__tmp["S"]
$__tmp = @{}
$__tmp["K"] = $ws2.Range("K8").Value2
$__tmp["L"] = $ws2.Range("L8").Value2
$__tmp["M"] = $ws2.Range("M8").Value2
$__tmp["N"] = $ws2.Range("N8").Value2
$__tmp["O"] = $ws2.Range("O8").Value2
$__tmp["P"] = $ws2.Range("P8").Value2
$__tmp["Q"] = $ws2.Range("Q8").Value2
$__tmp["R"] = $ws2.Range("R8").Value2
$__tmp["S"] = $ws2.Range("S8").Value2
$ws2.Range("N8").Value = $__tmp["K"]
$ws2.Range("O8").Value = $__tmp["L"]
$ws2.Range("P8").Value = $__tmp["M"]
$ws2.Range("Q8").Value = $__tmp["N"]
$ws2.Range("R8").Value = $__tmp["O"]
$ws2.Range("S8").Value = $__tmp["P"]
$ws2.Range("T8").Value = $__tmp["Q"]
$ws2.Range("U8").Value = $__tmp["R"]
$ws2.Range("V8").Value = $__tmp["S"]
$__tmp = @{}
$__tmp["K"] = $ws2.Range("K9").Value2
$__tmp["L"] = $ws2.Range("L9").Value2
$__tmp["M"] = $ws2.Range("M9").Value2
$__tmp["N"] = $ws2.Range("N9").Value2
$__tmp["O"] = $ws2.Range("O9").Value2
$__tmp["P"] = $ws2.Range("P9").Value2
$__tmp["Q"] = $ws2.Range("Q9").Value2
$__tmp["R"] = $ws2.Range("R9").Value2
$__tmp["S"] = $ws2.Range("S9").Value2
$ws2.Range("N9").Value = $__tmp["K"]
$ws2.Range("O9").Value = $__tmp["L"]
$ws2.Range("P9").Value = $__tmp["M"]
$ws2.Range("Q9").Value = $__tmp["N"]
$ws2.Range("R9").Value = $__tmp["O"]
$ws2.Range("S9").Value = $__tmp["P"]
$ws2.Range("T9").Value = $__tmp["Q"]
$ws2.Range("U9").Value = $__tmp["R"]
$ws2.Range("V9").Value = $__tmp["S"]
$__tmp = @{}
$__tmp["K"] = $ws2.Range("K10").Value2
$__tmp["L"] = $ws2.Range("L10").Value2
$__tmp["M"] = $ws2.Range("M10").Value2
$__tmp["N"] = $ws2.Range("N10").Value2
$__tmp["O"] = $ws2.Range("O10").Value2
$__tmp["P"] = $ws2.Range("P10").Value2
$__tmp["Q"] = $ws2.Range("Q10").Value2
$__tmp["R"] = $ws2.Range("R10").Value2
$__tmp["S"] = $ws2.Range("S10").Value2
$ws2.Range("N10").Value = $__tmp["K"]
$ws2.Range("O10").Value = $__tmp["L"]
$ws2.Range("P10").Value = $__tmp["M"]
$ws2.Range("Q10").Value = $__tmp["N"]
$ws2.Range("R10").Value = $__tmp["O"]
$ws2.Range("S10").Value = $__tmp["P"]
$ws2.Range("T10").Value = $__tmp["Q"]
$ws2.Range("U10").Value = $__tmp["R"]
$ws2.Range("V10").Value = $__tmp["S"]
$__tmp = @{}
$__tmp["K"] = $ws2.Range("K11").Value2
$__tmp["L"] = $ws2.Range("L11").Value2
$__tmp["M"] = $ws2.Range("M11").Value2
$__tmp["N"] = $ws2.Range("N11").Value2
$__tmp["O"] = $ws2.Range("O11").Value2
$__tmp["P"] = $ws2.Range("P11").Value2
$__tmp["Q"] = $ws2.Range("Q11").Value2
$__tmp["R"] = $ws2.Range("R11").Value2
$__tmp["S"] = $ws2.Range("S11").Value2
$ws2.Range("N11").Value = $__tmp["K"]
$ws2.Range("O11").Value = $__tmp["L"]
$ws2.Range("P11").Value = $__tmp["M"]
$ws2.Range("Q11").Value = $__tmp["N"]
$ws2.Range("R11").Value = $__tmp["O"]
$ws2.Range("S11").Value = $__tmp["P"]
$ws2.Range("T11").Value = $__tmp["Q"]
$ws2.Range("U11").Value = $__tmp["R"]
$ws2.Range("V11").Value = $__tmp["S"]
$__tmp = @{}
$__tmp["K"] = $ws2.Range("K12").Value2
$__tmp["L"] = $ws2.Range("L12").Value2
$__tmp["M"] = $ws2.Range("M12").Value2
$__tmp["N"] = $ws2.Range("N12").Value2
$__tmp["O"] = $ws2.Range("O12").Value2
$__tmp["P"] = $ws2.Range("P12").Value2
$__tmp["Q"] = $ws2.Range("Q12").Value2
$__tmp["R"] = $ws2.Range("R12").Value2
$__tmp["S"] = $ws2.Range("S12").Value2
$ws2.Range("N12").Value = $__tmp["K"]
$ws2.Range("O12").Value = $__tmp["L"]
$ws2.Range("P12").Value = $__tmp["M"]
$ws2.Range("Q12").Value = $__tmp["N"]
$ws2.Range("R12").Value = $__tmp["O"]
$ws2.Range("S12").Value = $__tmp["P"]
$ws2.Range("T12").Value = $__tmp["Q"]
$ws2.Range("U12").Value = $__tmp["R"]
$ws2.Range("V12").Value = $__tmp["S"]
$__tmp = @{}
$__tmp["K"] = $ws2.Range("K13").Value2
$__tmp["L"] = $ws2.Range("L13").Value2
$__tmp["M"] = $ws2.Range("M13").Value2
$__tmp["N"] = $ws2.Range("N13").Value2
$__tmp["O"] = $ws2.Range("O13").Value2
$__tmp["P"] = $ws2.Range("P13").Value2
$__tmp["Q"] = $ws2.Range("Q13").Value2
$__tmp["R"] = $ws2.Range("R13").Value2
$__tmp["S"] = $ws2.Range("S13").Value2
$ws2.Range("N13").Value = $__tmp["K"]
$ws2.Range("O13").Value = $__tmp["L"]
$ws2.Range("P13").Value = $__tmp["M"]
$ws2.Range("Q13").Value = $__tmp["N"]
$ws2.Range("R13").Value = $__tmp["O"]
$ws2.Range("S13").Value = $__tmp["P"]
$ws2.Range("T13").Value = $__tmp["Q"]
$ws2.Range("U13").Value = $__tmp["R"]
$ws2.Range("V13").Value = $__tmp["S"]
$__tmp = @{}
$__tmp["K"] = $ws2.Range("K14").Value2
$__tmp["L"] = $ws2.Range("L14").Value2
$__tmp["M"] = $ws2.Range("M14").Value2
$__tmp["N"] = $ws2.Range("N14").Value2
$__tmp["O"] = $ws2.Range("O14").Value2
$__tmp["P"] = $ws2.Range("P14").Value2
$__tmp["Q"] = $ws2.Range("Q14").Value2
$__tmp["R"] = $ws2.Range("R14").Value2
$__tmp["S"] = $ws2.Range("S14").Value2
$ws2.Range("N14").Value = $__tmp["K"]
$ws2.Range("O14").Value = $__tmp["L"]
$ws2.Range("P14").Value = $__tmp["M"]
$ws2.Range("Q14").Value = $__tmp["N"]
$ws2.Range("R14").Value = $__tmp["O"]
$ws2.Range("S14").Value = $__tmp["P"]
$ws2.Range("T14").Value = $__tmp["Q"]
$ws2.Range("U14").Value = $__tmp["R"]
$ws2.Range("V14").Value = $__tmp["S"]
$__tmp = @{}
$__tmp["K"] = $ws2.Range("K15").Value2
$__tmp["L"] = $ws2.Range("L15").Value2
$__tmp["M"] = $ws2.Range("M15").Value2
$__tmp["N"] = $ws2.Range("N15").Value2
$__tmp["O"] = $ws2.Range("O15").Value2
$__tmp["P"] = $ws2.Range("P15").Value2
$__tmp["Q"] = $ws2.Range("Q15").Value2
$__tmp["R"] = $ws2.Range("R15").Value2
$__tmp["S"] = $ws2.Range("S15").Value2
$ws2.Range("N15").Value = $__tmp["K"]
$ws2.Range("O15").Value = $__tmp["L"]
$ws2.Range("P15").Value = $__tmp["M"]
$ws2.Range("Q15").Value = $__tmp["N"]
$ws2.Range("R15").Value = $__tmp["O"]
$ws2.Range("S15").Value = $__tmp["P"]
$ws2.Range("T15").Value = $__tmp["Q"]
$ws2.Range("U15").Value = $__tmp["R"]
$ws2.Range("V15").Value = $__tmp["S"]
$__tmp = @{}
$__tmp["K"] = $ws2.Range("K16").Value2
$__tmp["L"] = $ws2.Range("L16").Value2
$__tmp["M"] = $ws2.Range("M16").Value2
$__tmp["N"] = $ws2.Range("N16").Value2
$__tmp["O"] = $ws2.Range("O16").Value2
$__tmp["P"] = $ws2.Range("P16").Value2
$__tmp["Q"] = $ws2.Range("Q16").Value2
$__tmp["R"] = $ws2.Range("R16").Value2
$__tmp["S"] = $ws2.Range("S16").Value2
$ws2.Range("N16").Value = $__tmp["K"]
$ws2.Range("O16").Value = $__tmp["L"]
$ws2.Range("P16").Value = $__tmp["M"]
$ws2.Range("Q16").Value = $__tmp["N"]
$ws2.Range("R16").Value = $__tmp["O"]
$ws2.Range("S16").Value = $__tmp["P"]
$ws2.Range("T16").Value = $__tmp["Q"]
$ws2.Range("U16").Value = $__tmp["R"]
$ws2.Range("V16").Value = $__tmp["S"]
$__tmp = @{}
$__tmp["K"] = $ws2.Range("K17").Value2
$__tmp["L"] = $ws2.Range("L17").Value2
$__tmp["M"] = $ws2.Range("M17").Value2
$__tmp["N"] = $ws2.Range("N17").Value2
$__tmp["O"] = $ws2.Range("O17").Value2
$__tmp["P"] = $ws2.Range("P17").Value2
$__tmp["Q"] = $ws2.Range("Q17").Value2
$__tmp["R"] = $ws2.Range("R17").Value2
$__tmp["S"] = $ws2.Range("S17").Value2
$ws2.Range("N17").Value = $__tmp["K"]
$ws2.Range("O17").Value = $__tmp["L"]
$ws2.Range("P17").Value = $__tmp["M"]
$ws2.Range("Q17").Value = $__tmp["N"]
$ws2.Range("R17").Value = $__tmp["O"]
$ws2.Range("S17").Value = $__tmp["P"]
$ws2.Range("T17").Value = $__tmp["Q"]
$ws2.Range("U17").Value = $__tmp["R"]
$ws2.Range("V17").Value = $__tmp["S"]
$__tmp = @{}
$__tmp["K"] = $ws2.Range("K18").Value2
$__tmp["L"] = $ws2.Range("L18").Value2
$__tmp["M"] = $ws2.Range("M18").Value2
$__tmp["N"] = $ws2.Range("N18").Value2
$__tmp["O"] = $ws2.Range("O18").Value2
$__tmp["P"] = $ws2.Range("P18").Value2
$__tmp["Q"] = $ws2.Range("Q18").Value2
$__tmp["R"] = $ws2.Range("R18").Value2
$__tmp["S"] = $ws2.Range("S18").Value2
$ws2.Range("N18").Value = $__tmp["K"]
$ws2.Range("O18").Value = $__tmp["L"]
$ws2.Range("P18").Value = $__tmp["M"]
$ws2.Range("Q18").Value = $__tmp["N"]
$ws2.Range("R18").Value = $__tmp["O"]
$ws2.Range("S18").Value = $__tmp["P"]
$ws2.Range("T18").Value = $__tmp["Q"]
$ws2.Range("U18").Value = $__tmp["R"]
$ws2.Range("V18").Value = $__tmp["S"]
$__tmp = @{}
$__tmp["K"] = $ws2.Range("K19").Value2
$__tmp["L"] = $ws2.Range("L19").Value2
$__tmp["M"] = $ws2.Range("M19").Value2
$__tmp["N"] = $ws2.Range("N19").Value2
$__tmp["O"] = $ws2.Range("O19").Value2
$__tmp["P"] = $ws2.Range("P19").Value2
$__tmp["Q"] = $ws2.Range("Q19").Value2
$__tmp["R"] = $ws2.Range("R19").Value2
$__tmp["S"] = $ws2.Range("S19").Value2
$ws2.Range("N19").Value = $__tmp["K"]
$ws2.Range("O19").Value = $__tmp["L"]
$ws2.Range("P19").Value = $__tmp["M"]
$ws2.Range("Q19").Value = $__tmp["N"]
$ws2.Range("R19").Value = $__tmp["O"]
$ws2.Range("S19").Value = $__tmp["P"]
$ws2.Range("T19").Value = $__tmp["Q"]
$ws2.Range("U19").Value = $__tmp["R"]
$ws2.Range("V19").Value = $__tmp["S"]
$__tmp = @{}
$__tmp["K"] = $ws2.Range("K20").Value2
$__tmp["L"] = $ws2.Range("L20").Value2
$__tmp["M"] = $ws2.Range("M20").Value2
$__tmp["N"] = $ws2.Range("N20").Value2
$__tmp["O"] = $ws2.Range("O20").Value2
$__tmp["P"] = $ws2.Range("P20").Value2
$__tmp["Q"] = $ws2.Range("Q20").Value2
$__tmp["R"] = $ws2.Range("R20").Value2
$__tmp["S"] = $ws2.Range("S20").Value2
$ws2.Range("N20").Value = $__tmp["K"]
$ws2.Range("O20").Value = $__tmp["L"]
$ws2.Range("P20").Value = $__tmp["M"]
$ws2.Range("Q20").Value = $__tmp["N"]
$ws2.Range("R20").Value = $__tmp["O"]
$ws2.Range("S20").Value = $__tmp["P"]
$ws2.Range("T20").Value = $__tmp["Q"]
$ws2.Range("U20").Value = $__tmp["R"]
$ws2.Range("V20").Value = $__tmp["S"]
$__tmp = @{}
$__tmp["K"] = $ws2.Range("K21").Value2
$__tmp["L"] = $ws2.Range("L21").Value2
$__tmp["M"] = $ws2.Range("M21").Value2
$__tmp["N"] = $ws2.Range("N21").Value2
$__tmp["O"] = $ws2.Range("O21").Value2
$__tmp["P"] = $ws2.Range("P21").Value2
$__tmp["Q"] = $ws2.Range("Q21").Value2
$__tmp["R"] = $ws2.Range("R21").Value2
$__tmp["S"] = $ws2.Range("S21").Value2
$ws2.Range("N21").Value = $__tmp["K"]
$ws2.Range("O21").Value = $__tmp["L"]
$ws2.Range("P21").Value = $__tmp["M"]
$ws2.Range("Q21").Value = $__tmp["N"]
$ws2.Range("R21").Value = $__tmp["O"]
$ws2.Range("S21").Value = $__tmp["P"]
$ws2.Range("T21").Value = $__tmp["Q"]
$ws2.Range("U21").Value = $__tmp["R"]
$ws2.Range("V21").Value = $__tmp["S"]
$__tmp = @{}
$__tmp["K"] = $ws2.Range("K22").Value2
$__tmp["L"] = $ws2.Range("L22").Value2
$__tmp["M"] = $ws2.Range("M22").Value2
$__tmp["N"] = $ws2.Range("N22").Value2
$__tmp["O"] = $ws2.Range("O22").Value2
$__tmp["P"] = $ws2.Range("P22").Value2
$__tmp["Q"] = $ws2.Range("Q22").Value2
$__tmp["R"] = $ws2.Range("R22").Value2
$__tmp["S"] = $ws2.Range("S22").Value2
$ws2.Range("N22").Value = $__tmp["K"]
$ws2.Range("O22").Value = $__tmp["L"]
$ws2.Range("P22").Value = $__tmp["M"]
$ws2.Range("Q22").Value = $__tmp["N"]
$ws2.Range("R22").Value = $__tmp["O"]
$ws2.Range("S22").Value = $__tmp["P"]
$ws2.Range("T22").Value = $__tmp["Q"]
$ws2.Range("U22").Value = $__tmp["R"]
$ws2.Range("V22").Value = $__tmp["S"]
$__tmp = @{}
$__tmp["K"] = $ws2.Range("K23").Value2
$__tmp["L"] = $ws2.Range("L23").Value2
$__tmp["M"] = $ws2.Range("M23").Value2
$__tmp["N"] = $ws2.Range("N23").Value2
$__tmp["O"] = $ws2.Range("O23").Value2
$__tmp["P"] = $ws2.Range("P23").Value2
$__tmp["Q"] = $ws2.Range("Q23").Value2
$__tmp["R"] = $ws2.Range("R23").Value2
$__tmp["S"] = $ws2.Range("S23").Value2
$ws2.Range("N23").Value = $__tmp["K"]
$ws2.Range("O23").Value = $__tmp["L"]
$ws2.Range("P23").Value = $__tmp["M"]
$ws2.Range("Q23").Value = $__tmp["N"]
$ws2.Range("R23").Value = $__tmp["O"]
$ws2.Range("S23").Value = $__tmp["P"]
$ws2.Range("T23").Value = $__tmp["Q"]
$ws2.Range("U23").Value = $__tmp["R"]
$ws2.Range("V23").Value = $__tmp["S"]
$__tmp = @{}
$__tmp["K"] = $ws2.Range("K24").Value2
$__tmp["L"] = $ws2.Range("L24").Value2
$__tmp["M"] = $ws2.Range("M24").Value2
$__tmp["N"] = $ws2.Range("N24").Value2
$__tmp["O"] = $ws2.Range("O24").Value2
$__tmp["P"] = $ws2.Range("P24").Value2
$__tmp["Q"] = $ws2.Range("Q24").Value2
$__tmp["R"] = $ws2.Range("R24").Value2
$__tmp["S"] = $ws2.Range("S24").Value2
$ws2.Range("N24").Value = $__tmp["K"]
$ws2.Range("O24").Value = $__tmp["L"]
$ws2.Range("P24").Value = $__tmp["M"]
$ws2.Range("Q24").Value = $__tmp["N"]
$ws2.Range("R24").Value = $__tmp["O"]
$ws2.Range("S24").Value = $__tmp["P"]
$ws2.Range("T24").Value = $__tmp["Q"]
$ws2.Range("U24").Value = $__tmp["R"]
$ws2.Range("V24").Value = $__tmp["S"]
$__tmp = @{}
$__tmp["K"] = $ws2.Range("K25").Value2
$__tmp["L"] = $ws2.Range("L25").Value2
$__tmp["M"] = $ws2.Range("M25").Value2
$__tmp["N"] = $ws2.Range("N25").Value2
$__tmp["O"] = $ws2.Range("O25").Value2
$__tmp["P"] = $ws2.Range("P25").Value2
$__tmp["Q"] = $ws2.Range("Q25").Value2
$__tmp["R"] = $ws2.Range("R25").Value2
$__tmp["S"] = $ws2.Range("S25").Value2
$ws2.Range("N25").Value = $__tmp["K"]
$ws2.Range("O25").Value = $__tmp["L"]
$ws2.Range("P25").Value = $__tmp["M"]
$ws2.Range("Q25").Value = $__tmp["N"]
$ws2.Range("R25").Value = $__tmp["O"]
$ws2.Range("S25").Value = $__tmp["P"]
$ws2.Range("T25").Value = $__tmp["Q"]
$ws2.Range("U25").Value = $__tmp["R"]
$ws2.Range("V25").Value = $__tmp["S"]
$__tmp = @{}
$__tmp["K"] = $ws2.Range("K26").Value2
$__tmp["L"] = $ws2.Range("L26").Value2
$__tmp["M"] = $ws2.Range("M26").Value2
$__tmp["N"] = $ws2.Range("N26").Value2
$__tmp["O"] = $ws2.Range("O26").Value2
$__tmp["P"] = $ws2.Range("P26").Value2
$__tmp["Q"] = $ws2.Range("Q26").Value2
$__tmp["R"] = $ws2.Range("R26").Value2
$__tmp["S"] = $ws2.Range("S26").Value2
$ws2.Range("N26").Value = $__tmp["K"]
$ws2.Range("O26").Value = $__tmp["L"]
$ws2.Range("P26").Value = $__tmp["M"]
$ws2.Range("Q26").Value = $__tmp["N"]
$ws2.Range("R26").Value = $__tmp["O"]
$ws2.Range("S26").Value = $__tmp["P"]
$ws2.Range("T26").Value = $__tmp["Q"]
$ws2.Range("U26").Value = $__tmp["R"]
$ws2.Range("V26").Value = $__tmp["S"]
$__tmp = @{}
$__tmp["K"] = $ws2.Range("K27").Value2
$__tmp["L"] = $ws2.Range("L27").Value2
$__tmp["M"] = $ws2.Range("M27").Value2
$__tmp["N"] = $ws2.Range("N27").Value2
$__tmp["O"] = $ws2.Range("O27").Value2
$__tmp["P"] = $ws2.Range("P27").Value2
$__tmp["Q"] = $ws2.Range("Q27").Value2
$__tmp["R"] = $ws2.Range("R27").Value2
$__tmp["S"] = $ws2.Range("S27").Value2
$ws2.Range("N27").Value = $__tmp["K"]
$ws2.Range("O27").Value = $__tmp["L"]
$ws2.Range("P27").Value = $__tmp["M"]
$ws2.Range("Q27").Value = $__tmp["N"]
$ws2.Range("R27").Value = $__tmp["O"]
$ws2.Range("S27").Value = $__tmp["P"]
$ws2.Range("T27").Value = $__tmp["Q"]
$ws2.Range("U27").Value = $__tmp["R"]
$ws2.Range("V27").Value = $__tmp["S"]
$__tmp = @{}
$__tmp["K"] = $ws2.Range("K28").Value2
$__tmp["L"] = $ws2.Range("L28").Value2
$__tmp["M"] = $ws2.Range("M28").Value2
$__tmp["N"] = $ws2.Range("N28").Value2
$__tmp["O"] = $ws2.Range("O28").Value2
$__tmp["P"] = $ws2.Range("P28").Value2
$__tmp["Q"] = $ws2.Range("Q28").Value2
$__tmp["R"] = $ws2.Range("R28").Value2
$__tmp["S"] = $ws2.Range("S28").Value2
$ws2.Range("N28").Value = $__tmp["K"]
$ws2.Range("O28").Value = $__tmp["L"]
$ws2.Range("P28").Value = $__tmp["M"]
$ws2.Range("Q28").Value = $__tmp["N"]
$ws2.Range("R28").Value = $__tmp["O"]
$ws2.Range("S28").Value = $__tmp["P"]
$ws2.Range("T28").Value = $__tmp["Q"]
$ws2.Range("U28").Value = $__tmp["R"]
$ws2.Range("V28").Value = $__tmp["S"]

# Sheet2 data: new temp_may/temp_jun/temp_jul columns (K, L, M)
$ws2.Range("K2").Value = 14.71047222
$ws2.Range("L2").Value = 16.400886159999999
$ws2.Range("M2").Value = 17.535275540000001
$ws2.Range("K3").Value = 14.597398310000001
$ws2.Range("L3").Value = 16.229580850000001
$ws2.Range("M3").Value = 17.35853943
$ws2.Range("K4").Value = 14.66313542
$ws2.Range("L4").Value = 16.28012996
$ws2.Range("M4").Value = 17.316293909999999
$ws2.Range("K5").Value = 14.80556548
$ws2.Range("L5").Value = 16.305771830000001
$ws2.Range("M5").Value = 17.47868369
$ws2.Range("K6").Value = 14.91295573
$ws2.Range("L6").Value = 16.460445679999999
$ws2.Range("M6").Value = 17.609603159999999
$ws2.Range("K7").Value = 15.02034598
$ws2.Range("L7").Value = 16.615119539999998
$ws2.Range("M7").Value = 17.740522630000001
$ws2.Range("K8").Value = 14.71770858
$ws2.Range("L8").Value = 16.304128720000001
$ws2.Range("M8").Value = 17.280460130000002
$ws2.Range("K9").Value = 14.60351389
$ws2.Range("L9").Value = 16.182056800000002
$ws2.Range("M9").Value = 17.17109722
$ws2.Range("K10").Value = 14.61614831
$ws2.Range("L10").Value = 16.234563489999999
$ws2.Range("M10").Value = 17.16211895
$ws2.Range("K11").Value = 13.996230779999999
$ws2.Range("L11").Value = 15.8404375
$ws2.Range("M11").Value = 16.883061380000001
$ws2.Range("K12").Value = 14.29384772
$ws2.Range("L12").Value = 16.16552927
$ws2.Range("M12").Value = 17.30660215
$ws2.Range("K13").Value = 14.14503925
$ws2.Range("L13").Value = 16.00298338
$ws2.Range("M13").Value = 17.094831769999999
$ws2.Range("K14").Value = 14.277400549999999
$ws2.Range("L14").Value = 16.10004588
$ws2.Range("M14").Value = 17.16383356
$ws2.Range("K15").Value = 14.380420389999999
$ws2.Range("L15").Value = 16.15130035
$ws2.Range("M15").Value = 17.224108650000002
$ws2.Range("K16").Value = 14.05140873
$ws2.Range("L16").Value = 15.83533482
$ws2.Range("M16").Value = 16.853850810000001
$ws2.Range("K17").Value = 14.918209450000001
$ws2.Range("L17").Value = 16.761449899999999
$ws2.Range("M17").Value = 17.767870739999999
$ws2.Range("K18").Value = 14.950884670000001
$ws2.Range("L18").Value = 16.722202630000002
$ws2.Range("M18").Value = 17.75468145
$ws2.Range("K19").Value = 14.64165017
$ws2.Range("L19").Value = 16.407492309999999
$ws2.Range("M19").Value = 17.37860573
$ws2.Range("K20").Value = 14.942062
$ws2.Range("L20").Value = 16.932037950000002
$ws2.Range("M20").Value = 18.06022536
$ws2.Range("K21").Value = 14.957031990000001
$ws2.Range("L21").Value = 16.90940179
$ws2.Range("M21").Value = 18.07080959
$ws2.Range("K22").Value = 14.82200645
$ws2.Range("L22").Value = 16.70975322
$ws2.Range("M22").Value = 17.86237388
$ws2.Range("K23").Value = 14.85719692
$ws2.Range("L23").Value = 17.066889880000002
$ws2.Range("M23").Value = 17.910047720000001
$ws2.Range("K24").Value = 15.00430283
$ws2.Range("L24").Value = 17.329374999999999
$ws2.Range("M24").Value = 18.222874099999999
$ws2.Range("K25").Value = 14.91703249
$ws2.Range("L25").Value = 17.146986859999998
$ws2.Range("M25").Value = 18.02783513
$ws2.Range("K26").Value = 15.06315551
$ws2.Range("L26").Value = 17.1429747
$ws2.Range("M26").Value = 18.040794129999998
$ws2.Range("K27").Value = 14.926790670000001
$ws2.Range("L27").Value = 16.989310759999999
$ws2.Range("M27").Value = 17.857511200000001
$ws2.Range("K28").Value = 14.92314807
$ws2.Range("L28").Value = 17.06262971
$ws2.Range("M28").Value = 18.0053573

# Sheet2 data: new fog_may/fog_jun/fog_jul columns (W, X, Y)
$ws2.Range("W2").Value = 0.10482180100000001
$ws2.Range("X2").Value = 0.082159624000000001
$ws2.Range("Y2").Value = 0.049212596999999997
$ws2.Range("W3").Value = 0.096114515999999997
$ws2.Range("X3").Value = 0.073459714999999995
$ws2.Range("Y3").Value = 0.039215687999999999
$ws2.Range("W4").Value = 0.096114515999999997
$ws2.Range("X4").Value = 0.073459714999999995
$ws2.Range("Y4").Value = 0.039215687999999999
$ws2.Range("W5").Value = 0.11273486200000001
$ws2.Range("X5").Value = 0.11032863699999999
$ws2.Range("Y5").Value = 0.04296875
$ws2.Range("W6").Value = 0.11273486200000001
$ws2.Range("X6").Value = 0.11032863699999999
$ws2.Range("Y6").Value = 0.04296875
$ws2.Range("W7").Value = 0.117154814
$ws2.Range("X7").Value = 0.11395348600000001
$ws2.Range("Y7").Value = 0.047058823999999999
$ws2.Range("W8").Value = 0.13953489099999999
$ws2.Range("X8").Value = 0.096698113000000002
$ws2.Range("Y8").Value = 0.045977010999999998
$ws2.Range("W9").Value = 0.16666667199999999
$ws2.Range("X9").Value = 0.110599078
$ws2.Range("Y9").Value = 0.055555555999999999
$ws2.Range("W10").Value = 0.16135881799999999
$ws2.Range("X10").Value = 0.11709602199999999
$ws2.Range("Y10").Value = 0.057361375999999999
$ws2.Range("W11").Value = 0.18670999999999999
$ws2.Range("X11").Value = 0.16153500000000001
$ws2.Range("Y11").Value = 0.083419999999999994
$ws2.Range("W12").Value = 0.154008433
$ws2.Range("X12").Value = 0.14841848599999999
$ws2.Range("Y12").Value = 0.079150579999999998
$ws2.Range("W13").Value = 0.18670999999999999
$ws2.Range("X13").Value = 0.16153500000000001
$ws2.Range("Y13").Value = 0.083419999999999994
$ws2.Range("W14").Value = 0.22698073099999999
$ws2.Range("X14").Value = 0.19002374999999999
$ws2.Range("Y14").Value = 0.10150376
$ws2.Range("W15").Value = 0.22698073099999999
$ws2.Range("X15").Value = 0.19002374999999999
$ws2.Range("Y15").Value = 0.10150376
$ws2.Range("W16").Value = 0.20212766500000001
$ws2.Range("X16").Value = 0.15813954199999999
$ws2.Range("Y16").Value = 0.073170729000000004
$ws2.Range("W17").Value = 0.25847458800000001
$ws2.Range("X17").Value = 0.20432692799999999
$ws2.Range("Y17").Value = 0.12098299
$ws2.Range("W18").Value = 0.26260504099999998
$ws2.Range("X18").Value = 0.20470587900000001
$ws2.Range("Y18").Value = 0.116412215
$ws2.Range("W19").Value = 0.241090149
$ws2.Range("X19").Value = 0.174712643
$ws2.Range("Y19").Value = 0.133079842
$ws2.Range("W20").Value = 0.25210085500000001
$ws2.Range("X20").Value = 0.22195121600000001
$ws2.Range("Y20").Value = 0.156673118
$ws2.Range("W21").Value = 0.25158563299999998
$ws2.Range("X21").Value = 0.22413793200000001
$ws2.Range("Y21").Value = 0.13035018700000001
$ws2.Range("W22").Value = 0.25158563299999998
$ws2.Range("X22").Value = 0.22413793200000001
$ws2.Range("Y22").Value = 0.13035018700000001
$ws2.Range("W23").Value = 0.284463882
$ws2.Range("X23").Value = 0.229468599
$ws2.Range("Y23").Value = 0.17307692799999999
$ws2.Range("W24").Value = 0.29787233499999999
$ws2.Range("X24").Value = 0.216152012
$ws2.Range("Y24").Value = 0.14779271199999999
$ws2.Range("W25").Value = 0.29094827200000001
$ws2.Range("X25").Value = 0.21980676099999999
$ws2.Range("Y25").Value = 0.152963668
$ws2.Range("W26").Value = 0.33916848900000002
$ws2.Range("X26").Value = 0.28605768100000001
$ws2.Range("Y26").Value = 0.194931775
$ws2.Range("W27").Value = 0.34341251900000003
$ws2.Range("X27").Value = 0.26763990500000001
$ws2.Range("Y27").Value = 0.18467582799999999
$ws2.Range("W28").Value = 0.34408602100000002
$ws2.Range("X28").Value = 0.26699030400000001
$ws2.Range("Y28").Value = 0.188349515

# Selections and active sheet
$ws1.Range("I11").Select()
$ws2.Activate()
$ws2.Range("M10").Select()

